$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking text values are stored as text (matching the
# original inlineStr cells) instead of being auto-converted to numbers/
# percentages by Excel. We set NumberFormat to Text ("@") on the D and E
# columns used below before writing the values.
$ws.Range("D2:E50").NumberFormat = "@"

# --- Updated price (D) / volume 1h (E) values ---
$ws.Range("D2").Value = "268.51"
$ws.Range("E2").Value = "2.60%"
$ws.Range("D3").Value = "26.66"
$ws.Range("E3").Value = "-2.02%"
$ws.Range("D4").Value = "4.705"
$ws.Range("E4").Value = "-0.07%"
$ws.Range("D5").Value = "0.06092"
$ws.Range("E5").Value = "-1.83%"
$ws.Range("D6").Value = "6.730"
$ws.Range("E6").Value = "0.21%"
$ws.Range("D7").Value = "0.8570"
$ws.Range("E7").Value = "0.86%"
$ws.Range("D8").Value = "0.8955"
$ws.Range("E8").Value = "-2.37%"
$ws.Range("E9").Value = "0.61%"
$ws.Range("D10").Value = "0.04920"
$ws.Range("E10").Value = "6.78%"
$ws.Range("D11").Value = "0.07082"
$ws.Range("E11").Value = "-0.02%"
$ws.Range("D12").Value = "0.03218"
$ws.Range("E12").Value = "2.11%"
$ws.Range("D13").Value = "0.09017"
$ws.Range("E13").Value = "-0.46%"
$ws.Range("E14").Value = "0.01%"
$ws.Range("D15").Value = "0.0006070"
$ws.Range("E15").Value = "-1.35%"
$ws.Range("D16").Value = "0.006026"
$ws.Range("E16").Value = "-0.29%"
$ws.Range("D17").Value = "3.461"
$ws.Range("E17").Value = "-0.19%"
$ws.Range("D18").Value = "3.165"
$ws.Range("E18").Value = "-0.07%"
$ws.Range("D19").Value = "2.242"
$ws.Range("E19").Value = "2.86%"
$ws.Range("D21").Value = "0.1298"
$ws.Range("E21").Value = "-0.79%"
$ws.Range("D22").Value = "3.846"
$ws.Range("E22").Value = "-6.09%"
$ws.Range("D23").Value = "0.04226"
$ws.Range("E23").Value = "-0.51%"
$ws.Range("D24").Value = "0.001184"
$ws.Range("E24").Value = "-2.16%"
$ws.Range("D25").Value = "0.004145"
$ws.Range("E25").Value = "9.03%"
$ws.Range("D26").Value = "0.0001200"
$ws.Range("E26").Value = "-0.05%"
$ws.Range("D27").Value = "0.0001679"
$ws.Range("E27").Value = "4.91%"
$ws.Range("D40").Value = "0.03941"
$ws.Range("E40").Value = "0.57%"
$ws.Range("E41").Value = "0.21%"
$ws.Range("D42").Value = "0.004184"
$ws.Range("E42").Value = "1.26%"
$ws.Range("D43").Value = "0.002010"
$ws.Range("E43").Value = "-7.95%"
$ws.Range("D44").Value = "0.01259"
$ws.Range("E44").Value = "-9.57%"
$ws.Range("D45").Value = "0.00005132"
$ws.Range("E45").Value = "-0.77%"
$ws.Range("E46").Value = "-0.05%"
$ws.Range("E49").Value = "-0.05%"
$ws.Range("E50").Value = "-0.05%"

# --- Rows 47 and 48 swapped (BOLO moved above CoinbaseStockToken) ---
# Row 47 becomes BOLO with refreshed price/volume data, row 48 becomes
# CoinbaseStockToken keeping its previous price/volume data.
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").Value = "0.4341"
$ws.Range("E47").Value = "158.90%"

$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "0.02448"
$ws.Range("E48").Value = "-31.81%"
